# Updated cryptos list on Wed Jun 26 15:26:59 UTC 2024 with GitHub Actions
#
# Refreshes the Price / Volume(1h) columns of the crypto table with new
# quotes, and swaps the ImmutableX / Aptos rows (rank 35 <-> 36).
#
# Price cells that look like plain numbers (single decimal point, e.g.
# "136.74") are forced to stay text via a temporary "@" (Text) number
# format - otherwise the COM layer would auto-coerce them to numeric
# values, which would change the cell's stored type. The format is reset
# back to "Normal" immediately afterwards so no residual cell-level style
# is left behind, matching the look of the untouched cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.448.79'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.361.48'
$ws.Range('E3').Value = '  -1.66%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').Value = '3.359.77'
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.74%  '
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.393'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').Value = '3.935.78'
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').Value = '3.347.64'
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('D18').Value = '61.509.27'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.98%  '
$ws.Range('E21').Value = '  -1.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '376.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('E23').Value = '  -3.89%  '
$ws.Range('D24').Value = '3.510.63'
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.70%  '
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('E31').Value = '  +3.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.80%  '
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.54'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.19'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.36%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.77'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.63%  '
$ws.Range('E40').Value = '  -4.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.771'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('E46').Value = '  -2.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.12'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.87'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('D50').Value = '2.370.22'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0261'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.71%  '
